$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("Complaint location", $true, $true, $false, $false, $false, $true, 1, $false, "Location/address", 2)
